# Update cryptocurrency price/volume data per upstream refresh (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.890.01"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.544.86"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'205.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -2.25%  "
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("D12").Value = "1.764.58"
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").Value = "1.548.22"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").Value = "'0.511"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("D16").Value = "26.864.76"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "'61.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "'213.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Value = "'7.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.34%  "
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").Value = "'4.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.68%  "
$ws.Range("D23").Value = "'9.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  -3.34%  "
$ws.Range("D25").Value = "'152.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("D27").Value = "'14.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("D33").Value = "1.357.17"
$ws.Range("E33").Value = "  -3.34%  "
$ws.Range("D34").Value = "'2.94"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'0.964"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.04%  "
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").Value = "'0.520"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.50%  "
$ws.Range("D40").Value = "'0.804"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("E42").Value = "  +3.52%  "
$ws.Range("D43").Value = "'0.988"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("E44").Value = "  +1.87%  "
$ws.Range("D45").Value = "'63.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("E46").Value = "  -2.07%  "
$ws.Range("D47").Value = "1.678.62"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("D48").Value = "'86.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").Value = "'0.0509"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("D50").Value = "0.0₇0970"
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("D51").Value = "'0.0948"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.02%  "
